$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('to believe', '信じる|しんじる'),
    @('confidence', '自信|じしん'),
    @('traffic signal', '信号|しんごう'),
    @('superstition', '迷信|めいしん'),
    @('experience', '経験|けいけん'),
    @('economy', '経済|けいざい'),
    @('oversensitive', '神経質|しんけいしつ'),
    @('management studies', '経営学|けいえいがく'),
    @('via Paris', 'パリ経由|パリけいゆ'),
    @('typhoon', '台風|たいふう'),
    @('two vehicles/machines/etc.', '二台|にだい'),
    @('stage', '舞台|ぶたい'),
    @('kitchen', '台所|だいどころ'),
    @('typhoon', '台風|たいふう'),
    @('bath', 'お風呂|おふろ'),
    @('Japanese style', '和風|わふう'),
    @('wind', '風|かぜ'),
    @('to catch a cold', '風邪をひく|かぜをひく'),
    @('dog', '犬|いぬ'),
    @('puppy', '子犬|こいぬ'),
    @('watch dog', '番犬|ばんけん'),
    @('seeing-eye dog', '盲導犬|もうどうけん'),
    @('heavy; serious (illness)', '重い|おもい'),
    @('to pile up', '重ねる|かさねる'),
    @('body weight', '体重|たいじゅう'),
    @('valuables', '貴重品|きちょうひん'),
    @('at first', '初めは|はじめは'),
    @('first time', '初めて|はじめて'),
    @('first', '最初|さいしょ'),
    @('first snow', '初雪|はつゆき'),
    @('first love', '初恋|はつこい'),
    @('young', '若い|わかい'),
    @('young people', '若者|わかもの'),
    @('youthful', '若々しい|わかわかしい'),
    @('to send', '送る|おくる'),
    @('sending money', '送金|そうきん'),
    @('postage', '送料|そうりょう'),
    @('broadcast', '放送|ほうそう'),
    @('out-of-service train', '回送電車|かいそうでんしゃ'),
    @('happy (しあわ)', '幸せな|しあわせな'),
    @('happy (こう)', '幸福な|こうふくな'),
    @('good fortune', '幸運|こううん'),
    @('misfortune', '不幸|ふこう'),
    @('fortunately', '幸い|さいわい'),
    @('watch; clock', '時計|とけい'),
    @('plan', '計画|けいかく'),
    @('sum', '合計|ごうけい'),
    @('to measure', '計る|はかる'),
    @('to be late (る)', '遅れる|おくれる'),
    @('to miss (transportation)', '乗り遅れる|のりおくれる'),
    @('late', '遅い|おそい'),
    @('to be late (する)', '遅刻する|ちこくする'),
    @('worried about', '心配な|しんぱいな'),
    @('to distribute', '配る|くばる'),
    @('delivery', '配達|はいたつ'),
    @('home delivery service', '宅配便|たくはいびん'),
    @('younger brother', '弟|おとうと'),
    @('brothers', '兄弟|きょうだい'),
    @('apprentice', '弟子|でし'),
    @('younger sister', '妹|いもうと'),
    @('sisters', '姉妹|しまい'),
    @('sister school', '姉妹校|しまいこう'),
)

$startRow = 58
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}
